$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "y_type"
$ws.Range("B2").Value = "c"
$ws.Range("A3").Value = "y_col"

$ws.Range("B4").Select()
